$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 43: Baby / Clean Bandit, Marina And The Diam / 11-11-2018 ---
$ws.Range("A43").Value = 41
$ws.Range("A42").Copy()
$ws.Range("A43").PasteSpecial(-4122)

$ws.Range("B43").Value = "Baby"
$ws.Range("C43").Value = "Clean Bandit, Marina And The Diam"

$ws.Range("D42").Copy()
$ws.Range("D43").PasteSpecial(-4104)

# --- Row 44: Bad Liar / Selena Gomez / 11-11-2018 ---
$ws.Range("A44").Value = 42
$ws.Range("A42").Copy()
$ws.Range("A44").PasteSpecial(-4122)

$ws.Range("B44").Value = "Bad Liar"
$ws.Range("C44").Value = "Selena Gomez"

$ws.Range("D42").Copy()
$ws.Range("D44").PasteSpecial(-4104)

$excel.CutCopyMode = 0

# Column B got visibly wider (target raw width 34.7109375 chars; the
# engine quantizes ColumnWidth to 1/6-character pixel steps, so 33.8 is
# the closest input that lands on the nearest achievable stored width).
$ws.Columns("B").ColumnWidth = 33.8
